$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on cells whose new numeric-looking values must remain as text
$textCells = "D5,D6,D9,D10,D11,D13,D15,D17,D18,D22,D23,D24,D26,D27,D28,D30,D31,D33,D35,D39,D40,D41,D42,D44,D45,D46,D47,D49".Split(",")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = '70.774.44'
$ws.Range("D3").Value = '3.611.97'
$ws.Range("E3").Value = '  +2.36%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '203.06'
$ws.Range("E5").Value = '  +3.95%  '
$ws.Range("D6").Value = '599.14'
$ws.Range("E6").Value = '  -1.29%  '
$ws.Range("E7").Value = '  +0.72%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = '0.217'
$ws.Range("E9").Value = '  +6.95%  '
$ws.Range("D10").Value = '0.646'
$ws.Range("E10").Value = '  +0.02%  '
$ws.Range("D11").Value = '53.98'
$ws.Range("E11").Value = '  +0.92%  '
$ws.Range("E12").Value = '  +0.61%  '
$ws.Range("D13").Value = '9.64'
$ws.Range("E13").Value = '  +1.90%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '4.184.84'
$ws.Range("E14").Value = '  +2.30%  '
$ws.Range("B15").Value = 'BitcoinCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D15").Value = '684.25'
$ws.Range("E15").Value = '  +15.51%  '
$ws.Range("D16").Value = '70.851.25'
$ws.Range("E16").Value = '  +1.38%  '
$ws.Range("D17").Value = '19.25'
$ws.Range("E17").Value = '  +1.45%  '
$ws.Range("D18").Value = '12.83'
$ws.Range("E18").Value = '  +0.87%  '
$ws.Range("D19").Value = '3.618.92'
$ws.Range("E19").Value = '  +2.70%  '
$ws.Range("E20").Value = '  +0.47%  '
$ws.Range("E21").Value = '  +1.81%  '
$ws.Range("D22").Value = '18.78'
$ws.Range("E22").Value = '  +5.45%  '
$ws.Range("D23").Value = '110.55'
$ws.Range("E23").Value = '  +7.53%  '
$ws.Range("D24").Value = '5.33'
$ws.Range("E24").Value = '  +3.29%  '
$ws.Range("E25").Value = '  -0.31%  '
$ws.Range("D26").Value = '3.04'
$ws.Range("E26").Value = '  -0.35%  '
$ws.Range("D27").Value = '10.64'
$ws.Range("E27").Value = '  -1.49%  '
$ws.Range("D28").Value = '6.01'
$ws.Range("E28").Value = '  -0.56%  '
$ws.Range("E29").Value = '  +6.53%  '
$ws.Range("D30").Value = '34.56'
$ws.Range("E30").Value = '  +4.15%  '
$ws.Range("D31").Value = '4.49'
$ws.Range("E31").Value = '  +6.26%  '
$ws.Range("E32").Value = '  +1.78%  '
$ws.Range("D33").Value = '12.29'
$ws.Range("E33").Value = '  -0.43%  '
$ws.Range("E34").Value = '  -0.20%  '
$ws.Range("D35").Value = '63.59'
$ws.Range("E35").Value = '  +0.28%  '
$ws.Range("E36").Value = '  +5.76%  '
$ws.Range("D37").Value = '3.873.37'
$ws.Range("E37").Value = '  +1.75%  '
$ws.Range("E38").Value = '  -0.11%  '
$ws.Range("D39").Value = '513.10'
$ws.Range("E39").Value = '  -0.26%  '
$ws.Range("D40").Value = '3.02'
$ws.Range("E40").Value = '  -5.57%  '
$ws.Range("D41").Value = '36.95'
$ws.Range("E41").Value = '  +1.22%  '
$ws.Range("D42").Value = '3.59'
$ws.Range("E42").Value = '  +0.27%  '
$ws.Range("E43").Value = '  -1.62%  '
$ws.Range("D44").Value = '0.138'
$ws.Range("E44").Value = '  +3.27%  '
$ws.Range("D45").Value = '0.0467'
$ws.Range("E45").Value = '  +4.33%  '
$ws.Range("D46").Value = '3.06'
$ws.Range("E46").Value = '  +8.66%  '
$ws.Range("D47").Value = '3.44'
$ws.Range("E47").Value = '  +5.32%  '
$ws.Range("E48").Value = '  +1.94%  '
$ws.Range("D49").Value = '8.65'
$ws.Range("E49").Value = '  +2.12%  '
$ws.Range("E50").Value = '  -0.25%  '
$ws.Range("E51").Value = '  +68.24%  '
